$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.273694276809692
$ws.Range("B1").Value = 2.257078886032104
$ws.Range("C1").Value = 4.45408821105957
$ws.Range("D1").Value = 2.93709135055542
$ws.Range("E1").Value = 1.012008786201477
